# Price-tracker update: append the 2026-02-07 scrape result as a new row
# at the bottom of the Date/Price/Discount/Incredible table.
#   Date=2026-02-07, Price=3199000, Discount=0, Incredible=0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the current table.
$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$rowRange = $ws.Range("A$newRow" + ":D$newRow")

# Every existing cell in this sheet stores its value as text (even the
# numeric-looking Price/Discount figures and the ISO dates), so force a
# text number format before writing the new values to avoid Excel
# auto-converting them into a real number / date serial.
$rowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = "2026-02-07"
$ws.Cells.Item($newRow, 2).Value = "3199000"
$ws.Cells.Item($newRow, 3).Value = "0"
$ws.Cells.Item($newRow, 4).Value = "0"

# Drop back to the workbook's default "Normal" style so the new row matches
# the unstyled look of every other row in the table.
$rowRange.Style = "Normal"
